$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("M2").Value = 1.743137
$ws.Range("N2").Value = 5.229411
$ws.Range("O2").Value = 0.03144673183548247
$ws.Range("P2").Value = 0.03144673183548247
$ws.Range("Q2").Value = 0.01979390168066667
$ws.Range("R2").Value = 0.178145115126
$ws.Range("S2").Value = 0.03144673183548247
$ws.Range("T2").Value = 0.03144673183548247

# Row 3
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("N3").Value = 3.848628
$ws.Range("O3").Value = 0.02314348071905789
$ws.Range("P3").Value = 0.02314348071905789
$ws.Range("Q3").Value = 0.01456748460533334
$ws.Range("R3").Value = 0.131107361448
$ws.Range("S3").Value = 0.02314348071905789
$ws.Range("T3").Value = 0.02314348071905789

# Row 4
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 52.405407
$ws.Range("N4").Value = 157.216221
$ws.Range("O4").Value = 0.9454097874454597
$ws.Range("P4").Value = 0.9454097874454597
$ws.Range("Q4").Value = 0.595080864954
$ws.Range("R4").Value = 5.355727784586001
$ws.Range("S4").Value = 0.9454097874454597
$ws.Range("T4").Value = 0.9454097874454597
